$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 (sensitivity_self_report): remove Harrison, 2019{#351} from study list; update primary_results
$ws.Range("B11").Value = '40 studies (Aita, 2018{#16}; Bakare, 2020{#45}; Bastiaens, 2017{#55}; Chiasson, 2012{#141}; Dakwar, 2012{#172}; Dunlop, 2018{#212}; Grogan, 2018{#323}; Houston, 2011{#390}; Kessler, 2005{#442}; Kwan, 2024{#490}; Lancaster, 2018{#495}; Lewandowski, 2008{#525}; Liu, 2023{#536}; Luty, 2009{#555}; Marchant, 2015{#571}; McCann, 2004{#603}; Mehringer, 2002{#611}; Palmer, 2023{#685}; Pettersson, 2018{#711}; Reimherr, 2021{#757}; Reyes, 2019{#762}; Singh, 2015{#857}; Solanto, 2004{#870}; Ustun, 2017{#940}; van de Glind, 2013{#942}; Van Voorhees, 2011{#949}; Vizgaitis, 2023{#960}; Young, 2023{#1016}; Young, 2016{#1017}; Kumar, 2011{#4185}; Roy-Byrne, 1997{#4485}; Erhardt, 1999{#5176}; Kingston, 2013{#5349}; De Quiros, 2001{#10619}; Dvorsky, 2016{#10626}; Faraone, 2010{#10633}; Skirrow, 2013{#10752}; Knight, 2025{#11044}; Herman, 2025{#14009}; Kitsune, 2016{#14746})'
$ws.Range("C11").Value = '100(#141, 490) to 28(#11044)'

# Row 12 (specificity_self_report): remove Harrison, 2019{#351} from study list
$ws.Range("B12").Value = '38 studies (Bakare, 2020{#45}; Bastiaens, 2017{#55}; Dakwar, 2012{#172}; Dunlop, 2018{#212}; Grogan, 2018{#323}; Houston, 2011{#390}; Kessler, 2005{#442}; Kwan, 2024{#490}; Lancaster, 2018{#495}; Lewandowski, 2008{#525}; Liu, 2023{#536}; Luty, 2009{#555}; Marchant, 2015{#571}; McCann, 2004{#603}; Mehringer, 2002{#611}; Palmer, 2023{#685}; Pettersson, 2018{#711}; Reimherr, 2021{#757}; Reyes, 2019{#762}; Singh, 2015{#857}; Solanto, 2004{#870}; Ustun, 2017{#940}; van de Glind, 2013{#942}; Van Voorhees, 2011{#949}; Vizgaitis, 2023{#960}; Young, 2023{#1016}; Young, 2016{#1017}; Kumar, 2011{#4185}; Roy-Byrne, 1997{#4485}; Erhardt, 1999{#5176}; Kingston, 2013{#5349}; De Quiros, 2001{#10619}; Dvorsky, 2016{#10626}; Faraone, 2010{#10633}; Skirrow, 2013{#10752}; Knight, 2025{#11044}; Herman, 2025{#14009}; Kitsune, 2016{#14746})'

# Row 83 (sensitivity_feigningADHD): remove Harrison, 2019{#351} from study list; update primary_results
$ws.Range("B83").Value = '22 studies (Abramson, 2023{#5}; Aita, 2018{#16}; Becke, 2023{#58}; Berger, 2021{#68}; Courrege, 2019{#164}; Musso, 2016{#646}; Potts, 2022{#733}; Quinn, 2003{#745}; Ramachandran, 2019{#750}; Robinson, 2023{#777}; Smith, 2017{#863}; Sollman, 2010{#872}; Spenceley, 2022{#883}; Young, 2011{#1015}; Harrison, 2020{#4042}; Morey, 2019{#4321}; Williamson, 2014{#4702}; Fuermaier, 2016{#5208}; Suhr, 2011{#10849}; Fuermaier, 2024{#11004}; Robinson, 2025{#11099}; Robinson, 2023{#14273})'
$ws.Range("C83").Value = '94(#745) to 24(#10849)'

# Row 84 (specificity_feigningADHD): remove Harrison, 2019{#351} from study list
$ws.Range("B84").Value = '20 studies (Abramson, 2023{#5}; Becke, 2023{#58}; Berger, 2021{#68}; Courrege, 2019{#164}; Musso, 2016{#646}; Potts, 2022{#733}; Quinn, 2003{#745}; Ramachandran, 2019{#750}; Robinson, 2023{#777}; Smith, 2017{#863}; Sollman, 2010{#872}; Spenceley, 2022{#883}; Young, 2011{#1015}; Harrison, 2020{#4042}; Morey, 2019{#4321}; Fuermaier, 2016{#5208}; Suhr, 2011{#10849}; Fuermaier, 2024{#11004}; Robinson, 2025{#11099}; Robinson, 2023{#14273})'

$wb.Save()
